$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G ("name") and give it the header "hhh".
# Everything that was in columns G:N shifts right into H:O.
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").Value = "hhh"

# Fix up the two data rows (2 and 3): amount formatting and the update-date string.
$ws.Range("A2:A3").Value = "300.000,99 MAD"
$ws.Range("N2:N3").Value = "Tue Sep 26 01:28:00 EDT 2023"
